$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASN")

# QTY column (H) was left blank for the line-item rows; fill it with qty 1
# as text, matching the rest of the sheet's inline-string cell formatting.
for ($r = 21; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = "1"
}
